$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the database/table/sequence setup tasks (rows 5-11, Status column F) as "Done"
$ws.Range("F5:F11").Value = "Done"

# Update the visible window / selection to match where the user was working
$ws.Range("D9").Select()
$excel.ActiveWindow.ScrollRow = 7
